# Auto-generated Excel COM-interop script to apply crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.57"
$ws.Range("E2").Value = "'0.64%"
$ws.Range("D3").Value = "'29.75"
$ws.Range("E3").Value = "'10.60%"
$ws.Range("D4").Value = "'5.174"
$ws.Range("E4").Value = "'2.10%"
$ws.Range("D5").Value = "'0.05710"
$ws.Range("D6").Value = "'6.602"
$ws.Range("E6").Value = "'1.82%"
$ws.Range("E7").Value = "'2.28%"
$ws.Range("D8").Value = "'0.8559"
$ws.Range("E8").Value = "'4.44%"
$ws.Range("D9").Value = "'0.8705"
$ws.Range("E9").Value = "'3.23%"
$ws.Range("D10").Value = "'0.1363"
$ws.Range("E10").Value = "'2.59%"
$ws.Range("D11").Value = "'0.07077"
$ws.Range("E11").Value = "'2.39%"
$ws.Range("D12").Value = "'0.02922"
$ws.Range("E12").Value = "'2.59%"
$ws.Range("D13").Value = "'0.09385"
$ws.Range("E13").Value = "'-0.05%"
$ws.Range("D14").Value = "'0.001512"
$ws.Range("E14").Value = "'-0.30%"
$ws.Range("D15").Value = "'0.04178"
$ws.Range("E15").Value = "'2.03%"
$ws.Range("D16").Value = "'0.0006021"
$ws.Range("E16").Value = "'-94.03%"
$ws.Range("D17").Value = "'0.006141"
$ws.Range("E17").Value = "'0.81%"
$ws.Range("E18").Value = "'3,767.27%"
$ws.Range("D19").Value = "'3.486"
$ws.Range("D20").Value = "'2.277"
$ws.Range("E20").Value = "'-1.70%"
$ws.Range("D22").Value = "'0.03365"
$ws.Range("E22").Value = "'5.31%"
$ws.Range("E23").Value = "'0.41%"
$ws.Range("D24").Value = "'3.468"
$ws.Range("E24").Value = "'-2.97%"
$ws.Range("E25").Value = "'0.48%"
$ws.Range("D26").Value = "'0.005026"
$ws.Range("E26").Value = "'26.75%"
$ws.Range("D27").Value = "'0.001220"
$ws.Range("E27").Value = "'0.20%"
$ws.Range("E28").Value = "'23.52%"
$ws.Range("D40").Value = "'0.03745"
$ws.Range("E40").Value = "'1.19%"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = "'0.005800"
$ws.Range("E41").Value = "'0.04%"
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1072"
$ws.Range("E42").Value = "'1.56%"
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.002000"
$ws.Range("E43").Value = "'-13.00%"
$ws.Range("D44").Value = "'0.009178"
$ws.Range("E44").Value = "'-2.36%"
$ws.Range("E45").Value = "'0.30%"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("D47").Value = "'0.05801"
$ws.Range("E47").Value = "'-51.64%"
$ws.Range("E48").Value = "'4.05%"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E50").Value = "'0.05%"
